# Re-salted output: the B column (Y-values feeding the scatter chart) was
# regenerated with a new random "salt" pass, and the active selection was
# left on M14 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New salted values for B1:B100 (paired with A1:A100 = 0..99, unchanged).
$newVals = @(
    28,21,-21,-28,-15,19,23,85,128,90,
    108,177,124,224,212,284,321,281,334,365,
    461,445,506,549,657,696,691,803,882,926,
    941,979,1063,1202,1263,1264,1335,1478,1565,1640,
    1725,1731,1838,1964,2071,2152,2245,2350,2450,2518,
    2579,2730,2775,2932,3017,3089,3265,3323,3525,3564,
    3690,3802,3994,4076,4239,4338,4444,4669,4794,4881,
    5072,5221,5318,5449,5615,5761,5920,6128,6292,6362,
    6586,6771,6935,7026,7268,7364,7540,7700,7895,8140,
    8296,8450,8680,8799,8991,9264,9422,9567,9784,10028
)

$count = $newVals.Length
$arr = New-Object 'object[,]' $count,1
for ($i = 0; $i -lt $count; $i++) {
    $arr[$i,0] = $newVals[$i]
}

$target = $ws.Range($ws.Cells.Item(1,2), $ws.Cells.Item($count,2))
$target.Value2 = $arr

# The workbook was left with this cell selected/active after editing.
$ws.Range("M14").Select()
